$wb = $excel.ActiveWorkbook

# ALC row 5
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 144.28572
$ws.Range("I5").Value = 134
$ws.Range("J5").Value = 170
$ws.Range("K5").Value = 134
$ws.Range("L5").Value = 170
$ws.Range("M5").Value = -19
$ws.Range("N5").Value = -400

# ALC row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 360.31818
$ws.Range("I6").Value = 15.166667
$ws.Range("K6").Value = 45.500001
$ws.Range("M6").Value = 66.499999

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1600.5
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 3402
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 3402
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -3540

# ALC row 95
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 15000
$ws.Range("J95").Value = 15000
$ws.Range("L95").Value = 15000
$ws.Range("N95").Value = -20492

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3249.75
$ws.Range("I106").Value = 2666.3333
$ws.Range("K106").Value = 2666.3333
$ws.Range("M106").Value = -2035.3333

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 547.6667
$ws.Range("I2").Value = 366.25
$ws.Range("K2").Value = 366.25
$ws.Range("M2").Value = -253.25

# ARM row 24
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 29000
$ws.Range("J24").Value = 29000
$ws.Range("L24").Value = 29000
$ws.Range("N24").Value = -29748

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12409.5
$ws.Range("I32").Value = 8768.223
$ws.Range("K32").Value = 8768.223
$ws.Range("M32").Value = -8481.223

# ARM row 100
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H100").Value = 29000
$ws.Range("J100").Value = 29000
$ws.Range("L100").Value = 29000
$ws.Range("N100").Value = -31164

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1161.1428
$ws.Range("I102").Value = 938
$ws.Range("K102").Value = 938
$ws.Range("M102").Value = 684

# ARM row 106
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 34757.75
$ws.Range("J106").Value = 34757.75
$ws.Range("L106").Value = 34757.75
$ws.Range("N106").Value = -37281.75

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 547.6667
$ws.Range("I116").Value = 366.25
$ws.Range("K116").Value = 366.25
$ws.Range("M116").Value = 1927.75

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 547.6667
$ws.Range("I3").Value = 366.25
$ws.Range("K3").Value = 366.25
$ws.Range("M3").Value = -252.25

# BSM row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 48664.832
$ws.Range("I82").Value = 12997.5
$ws.Range("J82").Value = 119999.5
$ws.Range("K82").Value = 12997.5
$ws.Range("L82").Value = 119999.5
$ws.Range("M82").Value = -12614.5
$ws.Range("N82").Value = -120765.5

# BSM row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 48664.832
$ws.Range("I85").Value = 12997.5
$ws.Range("J85").Value = 119999.5
$ws.Range("K85").Value = 12997.5
$ws.Range("L85").Value = 119999.5
$ws.Range("M85").Value = -11671.5
$ws.Range("N85").Value = -122651.5

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3464.1428
$ws.Range("I86").Value = 1812.25
$ws.Range("J86").Value = 5666.6665
$ws.Range("K86").Value = 1812.25
$ws.Range("L86").Value = 5666.6665
$ws.Range("M86").Value = -689.25
$ws.Range("N86").Value = -7912.6665

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3464.1428
$ws.Range("I89").Value = 1812.25
$ws.Range("J89").Value = 5666.6665
$ws.Range("K89").Value = 9061.25
$ws.Range("L89").Value = 28333.3325
$ws.Range("M89").Value = -3445.25
$ws.Range("N89").Value = -39565.3325

# BSM row 97
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 23350
$ws.Range("J97").Value = 50000
$ws.Range("L97").Value = 50000
$ws.Range("N97").Value = -51982

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 33.1
$ws.Range("I7").Value = 19.8
$ws.Range("K7").Value = 19.8
$ws.Range("M7").Value = 93.2

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3926.8
$ws.Range("I31").Value = 1645
$ws.Range("J31").Value = 5448
$ws.Range("K31").Value = 1645
$ws.Range("L31").Value = 5448
$ws.Range("M31").Value = -1350
$ws.Range("N31").Value = -6038

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3926.8
$ws.Range("I34").Value = 1645
$ws.Range("J34").Value = 5448
$ws.Range("K34").Value = 1645
$ws.Range("L34").Value = 5448
$ws.Range("M34").Value = -1443
$ws.Range("N34").Value = -5852

# CRP row 92
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# CRP row 95
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 27904.75
$ws.Range("J95").Value = 27904.75
$ws.Range("L95").Value = 27904.75
$ws.Range("N95").Value = -33396.75

# CRP row 96
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 16892.428
$ws.Range("J96").Value = 16892.428
$ws.Range("L96").Value = 16892.428
$ws.Range("N96").Value = -22384.428

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 4000
$ws.Range("I107").Value = 4000
$ws.Range("K107").Value = 4000
$ws.Range("M107").Value = -2080

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 316814.5
$ws.Range("J141").Value = 316814.5
$ws.Range("L141").Value = 316814.5
$ws.Range("N141").Value = -327174.5

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2052
$ws.Range("I139").Value = 2052
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 6156
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -1016
$ws.Range("N139").ClearContents()

# GSM row 101
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 7332.3335
$ws.Range("J101").Value = 7332.3335
$ws.Range("L101").Value = 7332.3335
$ws.Range("N101").Value = -13822.3335

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1782.8125
$ws.Range("J16").Value = 1499.8
$ws.Range("L16").Value = 1499.8
$ws.Range("N16").Value = -1839.8

# LTW row 32
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 938.5
$ws.Range("I32").Value = 1336
$ws.Range("J32").Value = 700
$ws.Range("K32").Value = 1336
$ws.Range("L32").Value = 700
$ws.Range("M32").Value = -1019
$ws.Range("N32").Value = -1334

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2167.5
$ws.Range("I100").Value = 2306.6667
$ws.Range("K100").Value = 2306.6667
$ws.Range("M100").Value = -1765.6667

# WVR row 58
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 41137.2
$ws.Range("I58").Value = 33833
$ws.Range("K58").Value = 33833
$ws.Range("M58").Value = -33525

# WVR row 92
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 32498.75
$ws.Range("J92").Value = 32498.75
$ws.Range("L92").Value = 32498.75
$ws.Range("N92").Value = -37490.75

# WVR row 104
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1958.625
$ws.Range("J136").Value = 3499.6667
$ws.Range("L136").Value = 10499.0001
$ws.Range("N136").Value = -15599.0001
